$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new row to the table (Table2); this extends the table range and
# worksheet dimension to A1:E23.
$tbl = $ws.ListObjects.Item(1)
$newRow = $tbl.ListRows.Add()

# Fill in the new row's values (row 23).
$ws.Cells.Item(23, 1).Value = "1729. Find Followers Count"
$ws.Cells.Item(23, 2).Value = "Easy"
$ws.Cells.Item(23, 2).Interior.Color = $ws.Cells.Item(22, 2).Interior.Color
$ws.Cells.Item(23, 3).Value = "Sorting and Grouping"

# Link cell with its hyperlink, matching the style used by the other Link cells.
$ws.Hyperlinks.Add($ws.Cells.Item(23, 5), "https://leetcode.com/problems/find-followers-count/solutions/1889791/mysql-ms-sql-oracle-simple-and-clean/?envType=study-plan-v2&envId=top-sql-50 ")
$ws.Cells.Item(23, 5).Style = "Hyperlink"

$ws.Cells.Item(23, 4).Value = "Use count(follower_id) as followers_count. Group by and order by user_id. Know that Group By performs the aggregation."

# Update the active selection/view like Excel would after editing the last row.
$ws.Activate() | Out-Null
$ws.Range("D24").Select() | Out-Null
